# Roma_matches_2020 - fill in xG/goals data for matches 10-15 (Napoli, Sassuolo,
# Bologna, Torino, Atalanta, Cagliari) that were previously blank in columns D-G.
#
# The source values are text-like numeric strings (matching the rest of the
# sheet, which stores xG_home/xG_away/goals_home/goals_away as shared-string
# text rather than numeric cells). A plain `.Value = "2.05867"` assignment
# gets auto-coerced to a real number by Excel, so instead we briefly stage the
# literal through a text formula and then flatten it back to a literal value
# via copy / paste-special(values) - this keeps the cell's stored type as
# text without touching any cell styles/number formats.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, [string]$text) {
    $cell = $ws.Range($range)
    $escaped = $text.Replace('"', '""')
    $cell.Formula = '="' + $escaped + '"'
    $cell.Copy()
    $cell.PasteSpecial(-4163)  # xlPasteValues
}

$newRows = @(
    @{ Row = 10; D = "2.05867";  E = "0.314633"; F = "4"; G = "0" },
    @{ Row = 11; D = "1.27993";  E = "0.209882"; F = "0"; G = "0" },
    @{ Row = 12; D = "0.354826"; E = "4.03543";  F = "1"; G = "5" },
    @{ Row = 13; D = "3.04581";  E = "0.975058"; F = "3"; G = "1" },
    @{ Row = 14; D = "1.91426";  E = "0.712758"; F = "4"; G = "1" },
    @{ Row = 15; D = "3.36165";  E = "1.19646";  F = "3"; G = "2" }
)

foreach ($entry in $newRows) {
    $r = $entry.Row
    Set-TextValue ("D" + $r) $entry.D
    Set-TextValue ("E" + $r) $entry.E
    Set-TextValue ("F" + $r) $entry.F
    Set-TextValue ("G" + $r) $entry.G
}

$excel.CutCopyMode = 0
